# Sept 2019 finance log - quick data cleaning; added missing negatives
#
# The rice-cooker purchase (row 28: Superstore / Cookware / =64.98*1.05) was
# returned, and the separate "Returned the rice cooker" row (row 29:
# Coffee Company / Coffee / -5.7) is actually an unrelated coffee expense
# that had accidentally been annotated as the return. Clean this up by:
#   1. Folding the refund into the original Superstore purchase row (row 7),
#      which now nets the rice cooker cost against the $68.23 refund and
#      carries the combined comment.
#   2. Deleting the now-redundant "Cookware" row (row 28) entirely, which
#      shifts every later row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (10-Sep Superstore/Food, -287.52): net out the returned rice cooker
# against its $68.23 refund, and update the comment to reflect the return.
$ws.Range("D7").Formula = "=-287.52+68.23"
$ws.Range("F7").Value = "Bough rice cooker, then returned later."

# Row 28 (16-Sep Superstore/Cookware, =64.98*1.05) is now redundant - remove
# it, shifting rows 29:52 up to 28:51.
$ws.Rows.Item(28).Delete()

# Match the author's final cursor position.
$ws.Range("F8").Select()
